$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 512
$ws.Range("I38").Value = 257.07693
$ws.Range("J38").Value = 1616.6666
$ws.Range("K38").Value = 771.2307900000001
$ws.Range("L38").Value = 4849.9998
$ws.Range("M38").Value = -399.2307900000001
$ws.Range("N38").Value = -5593.9998
# Row 40
$ws.Range("H40").Value = 4077.0789
$ws.Range("I40").Value = 4535.517
$ws.Range("J40").Value = 2599.889
$ws.Range("K40").Value = 4535.517
$ws.Range("L40").Value = 2599.889
$ws.Range("M40").Value = -4360.517
$ws.Range("N40").Value = -2949.889
# Row 47
$ws.Range("H47").Value = 23537
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 23537
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 23537
$ws.Range("M47").Value = ""
$ws.Range("N47").Value = -25481
# Row 87
$ws.Range("H87").Value = 29542.857
$ws.Range("J87").Value = 29542.857
$ws.Range("L87").Value = 29542.857
$ws.Range("N87").Value = -32038.857
# Row 90
$ws.Range("H90").Value = 29542.857
$ws.Range("J90").Value = 29542.857
$ws.Range("L90").Value = 88628.571
$ws.Range("N90").Value = -101108.571
# Row 113
$ws.Range("H113").Value = 3094.3845
$ws.Range("I113").Value = 2961
$ws.Range("K113").Value = 2961
$ws.Range("M113").Value = 293

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1710.5385
$ws.Range("I45").Value = 1793
$ws.Range("J45").Value = 1257
$ws.Range("K45").Value = 1793
$ws.Range("L45").Value = 1257
$ws.Range("M45").Value = -1416
$ws.Range("N45").Value = -2011
# Row 61
$ws.Range("H61").Value = 8404.143
$ws.Range("I61").Value = 8958.308000000001
$ws.Range("J61").Value = 1200
$ws.Range("K61").Value = 8958.308000000001
$ws.Range("L61").Value = 1200
$ws.Range("M61").Value = -8746.308000000001
$ws.Range("N61").Value = -1624
# Row 132
$ws.Range("H132").Value = 5305.5317
$ws.Range("I132").Value = 3486.7896
$ws.Range("J132").Value = 12984.667
$ws.Range("K132").Value = 10460.3688
$ws.Range("L132").Value = 38954.001
$ws.Range("M132").Value = -7930.3688
$ws.Range("N132").Value = -44014.001
# Row 136
$ws.Range("H136").Value = 8404.143
$ws.Range("I136").Value = 8958.308000000001
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 26874.924
$ws.Range("L136").Value = 3600
$ws.Range("M136").Value = -24324.924
$ws.Range("N136").Value = -8700

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1713.125
$ws.Range("I99").Value = 1808.4615
$ws.Range("J99").Value = 1300
$ws.Range("K99").Value = 1808.4615
$ws.Range("L99").Value = 1300
$ws.Range("M99").Value = -310.4614999999999
$ws.Range("N99").Value = -4296
# Row 107
$ws.Range("H107").Value = 1540.1111
$ws.Range("I107").Value = 1265.8572
$ws.Range("K107").Value = 1265.8572
$ws.Range("M107").Value = 654.1428000000001
# Row 137
$ws.Range("H137").Value = 71240
$ws.Range("J137").Value = 71240
$ws.Range("L137").Value = 71240
$ws.Range("N137").Value = -81440

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 94
$ws.Range("H94").Value = 803507.6
$ws.Range("I94").Value = 672174.7
$ws.Range("J94").Value = 1000507
$ws.Range("K94").Value = 672174.7
$ws.Range("L94").Value = 1000507
$ws.Range("M94").Value = -671723.7
$ws.Range("N94").Value = -1001409

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1745.7894
$ws.Range("I5").Value = 259.55554
$ws.Range("J5").Value = 2207.0344
$ws.Range("K5").Value = 778.66662
$ws.Range("L5").Value = 6621.1032
$ws.Range("M5").Value = -666.66662
$ws.Range("N5").Value = -6845.1032
# Row 36
$ws.Range("H36").Value = 6091.4
$ws.Range("I36").Value = 114.25
$ws.Range("K36").Value = 342.75
$ws.Range("M36").Value = -173.75
# Row 51
$ws.Range("H51").Value = 3001.4285
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3001.4285
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 9004.2855
$ws.Range("M51").Value = ""
$ws.Range("N51").Value = -9924.2855
# Row 68
$ws.Range("H68").Value = 17183.834
$ws.Range("I68").Value = 600
$ws.Range("J68").Value = 50351.5
$ws.Range("K68").Value = 1800
$ws.Range("L68").Value = 151054.5
$ws.Range("M68").Value = -989
$ws.Range("N68").Value = -152676.5
# Row 71
$ws.Range("H71").Value = 17183.834
$ws.Range("I71").Value = 600
$ws.Range("J71").Value = 50351.5
$ws.Range("K71").Value = 5400
$ws.Range("L71").Value = 453163.5
$ws.Range("M71").Value = -1344
$ws.Range("N71").Value = -461275.5
# Row 135
$ws.Range("H135").Value = 1745.7894
$ws.Range("I135").Value = 259.55554
$ws.Range("J135").Value = 2207.0344
$ws.Range("K135").Value = 2335.99986
$ws.Range("L135").Value = 19863.3096
$ws.Range("M135").Value = 199.0001400000001
$ws.Range("N135").Value = -24933.3096

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = ""
# Row 70
$ws.Range("H70").Value = 4269.087
$ws.Range("I70").Value = 4234
$ws.Range("J70").Value = 4328.9414
$ws.Range("K70").Value = 4234
$ws.Range("L70").Value = 4328.9414
$ws.Range("M70").Value = -3964
$ws.Range("N70").Value = -4868.9414
# Row 73
$ws.Range("H73").Value = 4269.087
$ws.Range("I73").Value = 4234
$ws.Range("J73").Value = 4328.9414
$ws.Range("K73").Value = 4234
$ws.Range("L73").Value = 4328.9414
$ws.Range("M73").Value = -3298
$ws.Range("N73").Value = -6200.9414
# Row 126
$ws.Range("H126").Value = 1882
$ws.Range("I126").Value = 1448
$ws.Range("J126").Value = 2750
$ws.Range("K126").Value = 4344
$ws.Range("L126").Value = 8250
$ws.Range("M126").Value = -1874
$ws.Range("N126").Value = -13190

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2188.3635
$ws.Range("I7").Value = 1872.5883
$ws.Range("J7").Value = 3262
$ws.Range("K7").Value = 1872.5883
$ws.Range("L7").Value = 3262
$ws.Range("M7").Value = -1760.5883
$ws.Range("N7").Value = -3486
# Row 100
$ws.Range("H100").Value = 2881.25
$ws.Range("I100").Value = 2457.1428
$ws.Range("K100").Value = 2457.1428
$ws.Range("M100").Value = -1916.1428
# Row 106
$ws.Range("H106").Value = 18000
$ws.Range("J106").Value = 18000
$ws.Range("L106").Value = 18000
$ws.Range("N106").Value = -20524
# Row 126
$ws.Range("H126").Value = 2188.3635
$ws.Range("I126").Value = 1872.5883
$ws.Range("J126").Value = 3262
$ws.Range("K126").Value = 5617.7649
$ws.Range("L126").Value = 9786
$ws.Range("M126").Value = -3147.7649
$ws.Range("N126").Value = -14726
# Row 132
$ws.Range("H132").Value = 9902.467000000001
$ws.Range("I132").Value = 14066.75
$ws.Range("J132").Value = 5143.2856
$ws.Range("K132").Value = 42200.25
$ws.Range("L132").Value = 15429.8568
$ws.Range("M132").Value = -39670.25
$ws.Range("N132").Value = -20489.8568
